$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (only changed cells per diff) ---
# Row 2
$ws.Range("T2").Value = 1.62

# Row 3
$ws.Range("L3").Value = 1.33
$ws.Range("T3").Value = 2.14
$ws.Range("AH3").Value = 38

# Row 4
$ws.Range("H4").Value = 3.9
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.01
$ws.Range("R4").Value = 1.19
$ws.Range("S4").Value = 2.72
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("V4").Value = 1.09
$ws.Range("W4").Value = 1.9
$ws.Range("X4").Value = 16
$ws.Range("Y4").Value = 22
$ws.Range("Z4").Value = 55
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 10
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 30
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 14.5
$ws.Range("AG4").Value = 14.5
$ws.Range("AH4").Value = 32
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 30
$ws.Range("AL4").Value = 65
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 22
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("I5").Value = 5.3
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 2.94
$ws.Range("O5").Value = 1.42
$ws.Range("R5").Value = 1.24
$ws.Range("S5").Value = 4.2
$ws.Range("T5").Value = 1.95
$ws.Range("U5").Value = 1.86
$ws.Range("V5").Value = 1.26
$ws.Range("W5").Value = 1.83
$ws.Range("X5").Value = 11.5
$ws.Range("Y5").Value = 14
$ws.Range("Z5").Value = 980
$ws.Range("AA5").Value = 130
$ws.Range("AB5").Value = 9.2
$ws.Range("AC5").Value = 9.2
$ws.Range("AD5").Value = 19
$ws.Range("AE5").Value = 70
$ws.Range("AF5").Value = 12.5
$ws.Range("AG5").Value = 11.5
$ws.Range("AH5").Value = 980
$ws.Range("AI5").Value = 85
$ws.Range("AJ5").Value = 28
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 980
$ws.Range("AM5").Value = 180
$ws.Range("AN5").Value = 980
$ws.Range("AO5").Value = 110

# --- Add new rows 6-11 ---
$bNewRows = $ws.Range("B6:B11")
$bNewRows.NumberFormat = "@"

# Row 6
$ws.Range("A6").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B6").Value = "2025-11-18"
$ws.Range("C6").Value = "22:00:00"
$ws.Range("D6").Value = "Panama"
$ws.Range("E6").Value = "El Salvador"
$ws.Range("F6").Value = 1.21
$ws.Range("G6").Value = 1.24
$ws.Range("H6").Value = 19.5
$ws.Range("I6").Value = 25
$ws.Range("J6").Value = 6.8
$ws.Range("K6").Value = 8.2
$ws.Range("L6").Value = 1.31
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 4.7
$ws.Range("O6").Value = 1.21
$ws.Range("P6").Value = 2.26
$ws.Range("Q6").Value = 1.62
$ws.Range("R6").Value = 1.51
$ws.Range("S6").Value = 2.58
$ws.Range("T6").Value = 2.54
$ws.Range("U6").Value = 1.52
$ws.Range("V6").Value = 1.04
$ws.Range("W6").Value = 5.1
$ws.Range("X6").Value = 980
$ws.Range("Y6").Value = 55
$ws.Range("Z6").Value = 290
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 8.6
$ws.Range("AC6").Value = 18.5
$ws.Range("AD6").Value = 85
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 7.2
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 60
$ws.Range("AI6").Value = 430
$ws.Range("AJ6").Value = 8.6
$ws.Range("AK6").Value = 17
$ws.Range("AL6").Value = 75
$ws.Range("AM6").Value = 450
$ws.Range("AN6").Value = 4.5
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("A7").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B7").Value = "2025-11-18"
$ws.Range("C7").Value = "22:00:00"
$ws.Range("D7").Value = "Haiti"
$ws.Range("E7").Value = "Nicaragua"
$ws.Range("F7").Value = 1.3
$ws.Range("G7").Value = 1.38
$ws.Range("H7").Value = 12
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 5.2
$ws.Range("K7").Value = 6.4
$ws.Range("L7").Value = 1.37
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 4.1
$ws.Range("O7").Value = 1.26
$ws.Range("P7").Value = 2.12
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.42
$ws.Range("S7").Value = 2.68
$ws.Range("T7").Value = 2.2
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 1.07
$ws.Range("W7").Value = 3.6
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 160
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 9.4
$ws.Range("AC7").Value = 15.5
$ws.Range("AD7").Value = 60
$ws.Range("AE7").Value = 310
$ws.Range("AF7").Value = 9
$ws.Range("AG7").Value = 11.5
$ws.Range("AH7").Value = 42
$ws.Range("AI7").Value = 240
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 19
$ws.Range("AL7").Value = 60
$ws.Range("AM7").Value = 270
$ws.Range("AN7").Value = 7
$ws.Range("AO7").Value = 510

# Row 8
$ws.Range("A8").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B8").Value = "2025-11-18"
$ws.Range("C8").Value = "22:00:00"
$ws.Range("D8").Value = "Costa Rica"
$ws.Range("E8").Value = "Honduras"
$ws.Range("F8").Value = 1.85
$ws.Range("G8").Value = 2.02
$ws.Range("H8").Value = 4.6
$ws.Range("I8").Value = 5.1
$ws.Range("J8").Value = 3.35
$ws.Range("K8").Value = 3.95
$ws.Range("L8").Value = 1.49
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.1
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 1.73
$ws.Range("Q8").Value = 2.12
$ws.Range("R8").Value = 1.26
$ws.Range("S8").Value = 3.9
$ws.Range("T8").Value = 1.95
$ws.Range("U8").Value = 1.86
$ws.Range("V8").Value = 1.25
$ws.Range("W8").Value = 1.99
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 15.5
$ws.Range("Z8").Value = 980
$ws.Range("AA8").Value = 160
$ws.Range("AB8").Value = 8
$ws.Range("AC8").Value = 8.6
$ws.Range("AD8").Value = 21
$ws.Range("AE8").Value = 90
$ws.Range("AF8").Value = 11.5
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 85
$ws.Range("AJ8").Value = 23
$ws.Range("AK8").Value = 24
$ws.Range("AL8").Value = 980
$ws.Range("AM8").Value = 180
$ws.Range("AN8").Value = 17.5
$ws.Range("AO8").Value = 130

# Row 9
$ws.Range("A9").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B9").Value = "2025-11-18"
$ws.Range("C9").Value = "22:00:00"
$ws.Range("D9").Value = "Jamaica"
$ws.Range("E9").Value = "Curacao"
$ws.Range("F9").Value = 2.08
$ws.Range("G9").Value = 2.28
$ws.Range("H9").Value = 3.95
$ws.Range("I9").Value = 4.5
$ws.Range("J9").Value = 3.05
$ws.Range("K9").Value = 3.45
$ws.Range("L9").Value = 1.52
$ws.Range("M9").Value = 1.09
$ws.Range("N9").Value = 2.82
$ws.Range("O9").Value = 1.45
$ws.Range("P9").Value = 1.62
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.23
$ws.Range("S9").Value = 4.4
$ws.Range("T9").Value = 1.97
$ws.Range("U9").Value = 1.84
$ws.Range("V9").Value = 1.29
$ws.Range("W9").Value = 1.78
$ws.Range("X9").Value = 12
$ws.Range("Y9").Value = 14.5
$ws.Range("Z9").Value = 34
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 9.2
$ws.Range("AC9").Value = 9
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 75
$ws.Range("AF9").Value = 15.5
$ws.Range("AG9").Value = 13.5
$ws.Range("AH9").Value = 25
$ws.Range("AI9").Value = 90
$ws.Range("AJ9").Value = 980
$ws.Range("AK9").Value = 980
$ws.Range("AL9").Value = 60
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 30
$ws.Range("AO9").Value = 1000

# Row 10
$ws.Range("A10").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B10").Value = "2025-11-18"
$ws.Range("C10").Value = "22:00:00"
$ws.Range("D10").Value = "Trinidad & Tobago"
$ws.Range("E10").Value = "Bermuda"
$ws.Range("F10").Value = 1.09
$ws.Range("G10").Value = 1.13
$ws.Range("H10").Value = 34
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 11.5
$ws.Range("K10").Value = 16
$ws.Range("L10").Value = 1.18
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 1.1
$ws.Range("O10").Value = 1.07
$ws.Range("P10").Value = 3.55
$ws.Range("Q10").Value = 1.27
$ws.Range("R10").Value = 2
$ws.Range("S10").Value = 1.61
$ws.Range("T10").Value = 2.24
$ws.Range("U10").Value = 1.38
$ws.Range("V10").Value = 1.02
$ws.Range("W10").Value = 8.6
$ws.Range("X10").Value = 70
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 21
$ws.Range("AC10").Value = 42
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 13
$ws.Range("AG10").Value = 25
$ws.Range("AH10").Value = 100
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 9.6
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 95
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000

# Row 11
$ws.Range("A11").Value = "FIFA World Cup Qualifiers - Americas"
$ws.Range("B11").Value = "2025-11-18"
$ws.Range("C11").Value = "22:00:00"
$ws.Range("D11").Value = "Guatemala"
$ws.Range("E11").Value = "Suriname"
$ws.Range("F11").Value = 2.76
$ws.Range("G11").Value = 3.4
$ws.Range("H11").Value = 2.54
$ws.Range("I11").Value = 2.9
$ws.Range("J11").Value = 3.2
$ws.Range("K11").Value = 3.7
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 1.63
$ws.Range("Q11").Value = 2.06
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("V11").Value = 0
$ws.Range("W11").Value = 0
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 0
$ws.Range("AD11").Value = 0
$ws.Range("AE11").Value = 0
$ws.Range("AF11").Value = 0
$ws.Range("AG11").Value = 0
$ws.Range("AH11").Value = 0
$ws.Range("AI11").Value = 0
$ws.Range("AJ11").Value = 0
$ws.Range("AK11").Value = 0
$ws.Range("AL11").Value = 0
$ws.Range("AM11").Value = 0
$ws.Range("AN11").Value = 0
$ws.Range("AO11").Value = 0

$bNewRows.ClearFormats()

